# Generate Report for Handback
#
# This script updates the localization-status workbook to reflect a
# handback event:
#   - Status text changes from "Ready for handoff" to
#     "Handed back: in sync with en-US" for both language sheets
#     (and, by extension, the Overview roll-up sheet).
#   - New "Latest Target File" (F) and "Latest Handback File" (G) columns
#     are populated (headers already exist) with hyperlinked file names,
#     for both the zh-cn and de-de language sheets.
#   - "Latest Handback DateTime" (H) is stamped with the handback time.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

$mdName425 = "425d36a0-c691-45e2-ab11-0a37ec306df5.md"
$mdNameA937 = "a9377a48-67ad-4e43-b924-f3c40a14ff5b.md"

$mdUrl425 = "https://github.com/OpenLocalizationTest/oltest/blob/c4dfe3db4d721c70c0e8f73746cc341293afba9a/e2e/425d36a0-c691-45e2-ab11-0a37ec306df5.md"
$mdUrlA937 = "https://github.com/OpenLocalizationTest/oltest/blob/c4dfe3db4d721c70c0e8f73746cc341293afba9a/e2e/a9377a48-67ad-4e43-b924-f3c40a14ff5b.md"

$xlfName425ZhCn = "425d36a0-c691-45e2-ab11-0a37ec306df5.30a09fc566298713cba6fd8bbaf67821415842ff.zh-cn.xlf"
$xlfNameA937ZhCn = "a9377a48-67ad-4e43-b924-f3c40a14ff5b.42150e9487e50f23fbeb2f83d65c9a386e6f9760.zh-cn.xlf"
$xlfUrl425ZhCn = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b091aaea3657432f9569a59562153bd30d8d7267/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/425d36a0-c691-45e2-ab11-0a37ec306df5.30a09fc566298713cba6fd8bbaf67821415842ff.zh-cn.xlf"
$xlfUrlA937ZhCn = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b091aaea3657432f9569a59562153bd30d8d7267/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a9377a48-67ad-4e43-b924-f3c40a14ff5b.42150e9487e50f23fbeb2f83d65c9a386e6f9760.zh-cn.xlf"

$xlfName425DeDe = "425d36a0-c691-45e2-ab11-0a37ec306df5.30a09fc566298713cba6fd8bbaf67821415842ff.de-de.xlf"
$xlfNameA937DeDe = "a9377a48-67ad-4e43-b924-f3c40a14ff5b.42150e9487e50f23fbeb2f83d65c9a386e6f9760.de-de.xlf"
$xlfUrl425DeDe = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a5197520729ec1994ad76e36cd2485f0ca487bef/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/425d36a0-c691-45e2-ab11-0a37ec306df5.30a09fc566298713cba6fd8bbaf67821415842ff.de-de.xlf"
$xlfUrlA937DeDe = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a5197520729ec1994ad76e36cd2485f0ca487bef/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a9377a48-67ad-4e43-b924-f3c40a14ff5b.42150e9487e50f23fbeb2f83d65c9a386e6f9760.de-de.xlf"

$handbackDateTimeZhCn = "2016-03-17 18:37:23"
$handbackDateTimeDeDe = "2016-03-17 18:37:29"

$linkColor = 15570276   # RGB(100,149,237) == #6495ED, matching the workbook's HyperLink style
$linkUnderline = 2      # xlUnderlineStyleSingle

function Style-AsHyperlink($range) {
    $range.Font.Underline = $linkUnderline
    $range.Font.Color = $linkColor
}

function Fill-LanguageSheet($sheetName, $xlfName425, $xlfNameA937, $xlfUrl425, $xlfUrlA937, $handbackDateTime) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Status column (C) -> handed back
    $ws.Range("C2").Value = $statusHandedBack
    $ws.Range("C3").Value = $statusHandedBack

    # Latest Target File (F) and Latest Handback File (G) for row 2 (425d36a0...)
    $ws.Hyperlinks.Add($ws.Range("F2"), $mdUrl425, "", "", $mdName425) | Out-Null
    Style-AsHyperlink $ws.Range("F2")

    $ws.Hyperlinks.Add($ws.Range("G2"), $xlfUrl425, "", "", $xlfName425) | Out-Null
    Style-AsHyperlink $ws.Range("G2")

    # Latest Target File (F) and Latest Handback File (G) for row 3 (a9377a48...)
    $ws.Hyperlinks.Add($ws.Range("F3"), $mdUrlA937, "", "", $mdNameA937) | Out-Null
    Style-AsHyperlink $ws.Range("F3")

    $ws.Hyperlinks.Add($ws.Range("G3"), $xlfUrlA937, "", "", $xlfNameA937) | Out-Null
    Style-AsHyperlink $ws.Range("G3")

    # Latest Handback DateTime (H)
    $ws.Range("H2").Value = $handbackDateTime
    $ws.Range("H3").Value = $handbackDateTime
}

Fill-LanguageSheet "zh-cn" $xlfName425ZhCn $xlfNameA937ZhCn $xlfUrl425ZhCn $xlfUrlA937ZhCn $handbackDateTimeZhCn
Fill-LanguageSheet "de-de" $xlfName425DeDe $xlfNameA937DeDe $xlfUrl425DeDe $xlfUrlA937DeDe $handbackDateTimeDeDe

# The Overview sheet's zh-cn/de-de status columns shared the same string
# as the language sheets' Status column; keep them in sync explicitly.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusHandedBack
$wsOverview.Range("C2").Value = $statusHandedBack
$wsOverview.Range("B3").Value = $statusHandedBack
$wsOverview.Range("C3").Value = $statusHandedBack

Write-Host "Handback report generated."
